# Apply "fixed workflow" re-run: both sheets (NBR, BAR) now reflect cutoffs
# 4-7 minutes later in the sweep (B column 5..19 instead of 1..19), which
# drops the last 4 rows of each sheet's 19-row table down to 15 rows
# (A1:C16 instead of A1:C20). Column A (the 0-based cutoff index) is left
# untouched; only columns B and C are rewritten, then the now-unused
# trailing rows are removed.

$wb = $excel.ActiveWorkbook

$newB = @(5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19)

$sheetData = @{
    "NBR" = @(806, 812, 804, 804, 794, 771, 790, 782, 758, 766, 757, 753, 793, 777, 740)
    "BAR" = @(966, 971, 960, 957, 959, 953, 924, 916, 918, 903, 903, 897, 854, 858, 888)
}

foreach ($ws in $wb.Worksheets) {
    $newC = $sheetData[$ws.Name]
    if ($newC -eq $null) { continue }

    for ($i = 0; $i -lt $newB.Length; $i++) {
        $row = 2 + $i
        $ws.Cells.Item($row, 2).Value = $newB[$i]
        $ws.Cells.Item($row, 3).Value = $newC[$i]
    }

    $ws.Rows("17:20").Delete()
}
